# Auto update Excel log
# Appends new sensor-log rows to the ALERTS, PIR, Humidity and Proximity
# sheets, matching the latest capture batch (2026-01-30, ~17:52-17:55).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALERTS sheet: rows 18-20 (row number, Date, Timestamp, Hour, Location,
# Value, Status)
# ---------------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")
$wsAlerts.Range("A18:F20").NumberFormat = "@"

$alertsData = @(
    ,@(18,"2026-01-30","17:54:08","17:00","Living Room","CRITICAL","FALL_DETECTED")
    ,@(19,"2026-01-30","17:54:12","17:00","Living Room","CRITICAL","FALL_DETECTED")
    ,@(20,"2026-01-30","17:54:58","17:00","Bathroom","MINIMAL","MINIMAL ALERT: Bathroom occupied, no motion > 20s.")
)

foreach ($row in $alertsData) {
    $r = $row[0]
    $wsAlerts.Cells.Item($r, 1).Value = $row[1]
    $wsAlerts.Cells.Item($r, 2).Value = $row[2]
    $wsAlerts.Cells.Item($r, 3).Value = $row[3]
    $wsAlerts.Cells.Item($r, 4).Value = $row[4]
    $wsAlerts.Cells.Item($r, 5).Value = $row[5]
    $wsAlerts.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# PIR sheet: rows 370-388
# ---------------------------------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")
$wsPir.Range("A370:F388").NumberFormat = "@"

$pirData = @(
    ,@(370,"2026-01-30","17:52:55","17:00","Bathroom","No Motion","Inactive")
    ,@(371,"2026-01-30","17:52:57","17:00","Bathroom","No Motion","Inactive")
    ,@(372,"2026-01-30","17:53:02","17:00","Bathroom","No Motion","Inactive")
    ,@(373,"2026-01-30","17:53:07","17:00","Bathroom","No Motion","Inactive")
    ,@(374,"2026-01-30","17:53:12","17:00","Bathroom","No Motion","Inactive")
    ,@(375,"2026-01-30","17:53:17","17:00","Bathroom","No Motion","Inactive")
    ,@(376,"2026-01-30","17:53:22","17:00","Bathroom","No Motion","Inactive")
    ,@(377,"2026-01-30","17:54:13","17:00","Bathroom","No Motion","Inactive")
    ,@(378,"2026-01-30","17:54:13","17:00","Bathroom","No Motion","Inactive")
    ,@(379,"2026-01-30","17:54:17","17:00","Bathroom","No Motion","Inactive")
    ,@(380,"2026-01-30","17:54:22","17:00","Bathroom","No Motion","Inactive")
    ,@(381,"2026-01-30","17:54:27","17:00","Bathroom","No Motion","Inactive")
    ,@(382,"2026-01-30","17:54:32","17:00","Bathroom","No Motion","Inactive")
    ,@(383,"2026-01-30","17:54:37","17:00","Bathroom","No Motion","Inactive")
    ,@(384,"2026-01-30","17:54:42","17:00","Bathroom","No Motion","Inactive")
    ,@(385,"2026-01-30","17:54:47","17:00","Bathroom","No Motion","Inactive")
    ,@(386,"2026-01-30","17:54:52","17:00","Bathroom","No Motion","Inactive")
    ,@(387,"2026-01-30","17:54:57","17:00","Bathroom","No Motion","Inactive")
    ,@(388,"2026-01-30","17:55:02","17:00","Bathroom","No Motion","Inactive")
)

foreach ($row in $pirData) {
    $r = $row[0]
    $wsPir.Cells.Item($r, 1).Value = $row[1]
    $wsPir.Cells.Item($r, 2).Value = $row[2]
    $wsPir.Cells.Item($r, 3).Value = $row[3]
    $wsPir.Cells.Item($r, 4).Value = $row[4]
    $wsPir.Cells.Item($r, 5).Value = $row[5]
    $wsPir.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Humidity sheet: rows 255-263
# ---------------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$wsHumidity.Range("A255:F263").NumberFormat = "@"

$humidityData = @(
    ,@(255,"2026-01-30","17:52:58","17:00","Bathroom","85.9%","Active")
    ,@(256,"2026-01-30","17:53:03","17:00","Bathroom","86.8%","Active")
    ,@(257,"2026-01-30","17:53:07","17:00","Bathroom","85.8%","Active")
    ,@(258,"2026-01-30","17:53:23","17:00","Bathroom","86.8%","Active")
    ,@(259,"2026-01-30","17:54:18","17:00","Bathroom","86.7%","Active")
    ,@(260,"2026-01-30","17:54:43","17:00","Bathroom","86.7%","Active")
    ,@(261,"2026-01-30","17:54:48","17:00","Bathroom","86.8%","Active")
    ,@(262,"2026-01-30","17:54:59","17:00","Bathroom","86.8%","Active")
    ,@(263,"2026-01-30","17:55:03","17:00","Bathroom","86.8%","Active")
)

foreach ($row in $humidityData) {
    $r = $row[0]
    $wsHumidity.Cells.Item($r, 1).Value = $row[1]
    $wsHumidity.Cells.Item($r, 2).Value = $row[2]
    $wsHumidity.Cells.Item($r, 3).Value = $row[3]
    $wsHumidity.Cells.Item($r, 4).Value = $row[4]
    $wsHumidity.Cells.Item($r, 5).Value = $row[5]
    $wsHumidity.Cells.Item($r, 6).Value = $row[6]
}

# ---------------------------------------------------------------------------
# Proximity sheet: rows 60-64
# ---------------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
$wsProximity.Range("A60:F64").NumberFormat = "@"

$proximityData = @(
    ,@(60,"2026-01-30","17:52:56","17:00","Bathroom Door","EXIT","User EXITED Bathroom")
    ,@(61,"2026-01-30","17:53:05","17:00","Bathroom Door","ENTER","User ENTERED Bathroom")
    ,@(62,"2026-01-30","17:53:15","17:00","Bathroom Door","EXIT","User EXITED Bathroom")
    ,@(63,"2026-01-30","17:53:23","17:00","Bathroom Door","ENTER","User ENTERED Bathroom")
    ,@(64,"2026-01-30","17:54:34","17:00","Bathroom Door","ENTER","User ENTERED Bathroom")
)

foreach ($row in $proximityData) {
    $r = $row[0]
    $wsProximity.Cells.Item($r, 1).Value = $row[1]
    $wsProximity.Cells.Item($r, 2).Value = $row[2]
    $wsProximity.Cells.Item($r, 3).Value = $row[3]
    $wsProximity.Cells.Item($r, 4).Value = $row[4]
    $wsProximity.Cells.Item($r, 5).Value = $row[5]
    $wsProximity.Cells.Item($r, 6).Value = $row[6]
}
